$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.961.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.104.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.098.29"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.75%  "
$ws.Range("E9").Value = "  +1.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +13.05%  "
$ws.Range("E11").Value = "  +2.54%  "
$ws.Range("E12").Value = "  +3.67%  "
$ws.Range("E13").Value = "  +4.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.606.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.048.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.103.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.61%  "
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "481.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.74%  "
$ws.Range("E21").Value = "  +3.43%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("E23").Value = "  +6.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +11.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  +2.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.34%  "
$ws.Range("E29").Value = "  +5.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("E32").Value = "  +2.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.67"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.35%  "
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "473.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.07%  "
$ws.Range("E38").Value = "  +5.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0831"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +18.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.017.94"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.52%  "
$ws.Range("E42").Value = "  +1.22%  "
$ws.Range("E43").Value = "  -2.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "28.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.83%  "
$ws.Range("E45").Value = "  +5.18%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.66%  "
$ws.Range("E48").Value = "  +2.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₃0523"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "116.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.50%  "
